$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.99999998969518844
$ws.Range("A2").Value = 0.99924244051094935
$ws.Range("A3").Value = 0.99844862777366206
$ws.Range("A4").Value = 1.0037340671260342
$ws.Range("A5").Value = 0.99595523261256202
$ws.Range("A6").Value = 0.97948246301338859
$ws.Range("A7").Value = 0.97808956863881791
$ws.Range("A8").Value = 0.97570164216292243
$ws.Range("A9").Value = 0.9760879567701668
$ws.Range("A10").Value = 0.97742931276151235
$ws.Range("A11").Value = 0.97766999223228845
$ws.Range("A12").Value = 0.97838523294538193
$ws.Range("A13").Value = 0.96541393570186274
$ws.Range("A14").Value = 0.96155575153418726
$ws.Range("A15").Value = 0.95851736445663538
$ws.Range("A16").Value = 0.9550823008375906
$ws.Range("A17").Value = 0.95137481787046818
$ws.Range("A18").Value = 0.95026594657892449
$ws.Range("A19").Value = 0.99840474237769306
$ws.Range("A20").Value = 0.99313090178385699
$ws.Range("A21").Value = 0.99173246876846932
$ws.Range("A22").Value = 0.99046797959267407
$ws.Range("A23").Value = 0.99141488981521886
$ws.Range("A24").Value = 0.97839508059360836
$ws.Range("A25").Value = 0.97193827748871664
$ws.Range("A26").Value = 0.97704693499903517
$ws.Range("A27").Value = 0.97390020125230936
$ws.Range("A28").Value = 0.96285052984369146
$ws.Range("A29").Value = 0.95530379202189386
$ws.Range("A30").Value = 0.95164712637132087
$ws.Range("A31").Value = 0.95447376803280914
$ws.Range("A32").Value = 0.9534647350178167
$ws.Range("A33").Value = 0.95294473926610257
